{"js": "// Update the date paragraph and the 25 multiplication-problem table cells\n// to the new values from the commit diff. Cells are addressed by their\n// (row, column) position in the table so that values which coincidentally\n// collide with other old/new values (e.g. \"14\u00d758=812\" appears both as an\n// old value in row 0 and as a new value in row 3) are never mixed up by a\n// text search.\n\n// 1) Update the date line above the table.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\n// Guard against re-running on an already-updated document, but always land\n// on the target value either way.\nif (dateParagraph.text !== \"2025-05-09 Friday\") {\n  dateParagraph.getRange().insertText(\"2025-05-09 Friday\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Update the 25 multiplication cells, addressed by (row, col).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row indices (0-based) of the table rows that actually hold the visible\n// multiplication problems; the rows in between are blank spacer rows.\nconst dataRowIndices = [0, 4, 9, 14, 19];\n\nconst newValues = [\n  [\"41\u00d730=1230\", \"77\u00d799=7623\", \"94\u00d715=1410\", \"19\u00d760=1140\", \"98\u00d775=7350\"],\n  [\"71\u00d785=6035\", \"45\u00d781=3645\", \"62\u00d778=4836\", \"48\u00d752=2496\", \"50\u00d764=3200\"],\n  [\"51\u00d731=1581\", \"34\u00d753=1802\", \"83\u00d715=1245\", \"81\u00d752=4212\", \"78\u00d772=5616\"],\n  [\"16\u00d798=1568\", \"14\u00d758=812\",  \"80\u00d730=2400\", \"73\u00d735=2555\", \"75\u00d732=2400\"],\n  [\"98\u00d757=5586\", \"99\u00d746=4554\", \"52\u00d759=3068\", \"74\u00d748=3552\", \"32\u00d732=1024\"],\n];\n\nfor (let r = 0; r < dataRowIndices.length; r++) {\n  const rowIndex = dataRowIndices[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(rowIndex, c);\n    cell.value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 multiplication-problem table cells to the\n# new values from the commit diff. Cells are addressed by their (row, column)\n# position in the table (Word COM uses 1-based indices) so that values which\n# coincidentally collide with other old/new values (e.g. \"14\u00d758=812\" is an\n# old value in row 1 and also a new value in row 4) are never mixed up by a\n# text search/replace.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph above the table (guard against re-running on\n#    an already-updated document, but always land on the target value).\n$dateRange = $d.Paragraphs.Item(1).Range\nif ($dateRange.Text.TrimEnd(\"`r\", \"`a\") -ne \"2025-05-09 Friday\") {\n    $d.Paragraphs.Item(1).Range.Text = \"2025-05-09 Friday\"\n}\n\n# 2) Update the 25 multiplication cells, addressed by (row, col), 1-based.\n$table = $d.Tables.Item(1)\n\n# Row indices (1-based) of the table rows that actually hold the visible\n# multiplication problems; the rows in between are blank spacer rows.\n$dataRowIndices = @(1, 5, 10, 15, 20)\n\n$newValues = @(\n    @(\"41\u00d730=1230\", \"77\u00d799=7623\", \"94\u00d715=1410\", \"19\u00d760=1140\", \"98\u00d775=7350\"),\n    @(\"71\u00d785=6035\", \"45\u00d781=3645\", \"62\u00d778=4836\", \"48\u00d752=2496\", \"50\u00d764=3200\"),\n    @(\"51\u00d731=1581\", \"34\u00d753=1802\", \"83\u00d715=1245\", \"81\u00d752=4212\", \"78\u00d772=5616\"),\n    @(\"16\u00d798=1568\", \"14\u00d758=812\",  \"80\u00d730=2400\", \"73\u00d735=2555\", \"75\u00d732=2400\"),\n    @(\"98\u00d757=5586\", \"99\u00d746=4554\", \"52\u00d759=3068\", \"74\u00d748=3552\", \"32\u00d732=1024\")\n)\n\nfor ($r = 0; $r -lt $dataRowIndices.Length; $r++) {\n    $rowIndex = $dataRowIndices[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $table.Cell($rowIndex, $c).Range.Text = $newValues[$r][$c - 1]\n    }\n}\n"}
